$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 16,14
$arr[0,0] = 20.002957
$arr[0,1] = 60.008871
$arr[0,2] = 0.7920860939997775
$arr[0,3] = 0.7920860939997775
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 3.151158666666667
$arr[0,7] = 9.453476
$arr[0,8] = 0.03114707555614071
$arr[0,9] = 0.03114707555614071
$arr[0,10] = 63.03249130951067
$arr[0,11] = 567.292421785596
$arr[0,12] = 0.02467116541677944
$arr[0,13] = 0.02467116541677944
$arr[1,0] = 20.002957
$arr[1,1] = 60.008871
$arr[1,2] = 0.7920860939997775
$arr[1,3] = 0.7920860939997775
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 5.038243666666667
$arr[1,7] = 15.114731
$arr[1,8] = 0.04979963650066307
$arr[1,9] = 0.04979963650066306
$arr[1,10] = 100.7797714198557
$arr[1,11] = 907.0179427787011
$arr[1,12] = 0.03944559955841895
$arr[1,13] = 0.03944559955841895
$arr[2,0] = 20.002957
$arr[2,1] = 60.008871
$arr[2,2] = 0.7920860939997775
$arr[2,3] = 0.7920860939997775
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 92.91163899999999
$arr[2,7] = 278.734917
$arr[2,8] = 0.9183688116343246
$arr[2,9] = 0.9183688116343246
$arr[2,10] = 1858.507519716523
$arr[2,11] = 16726.56767744871
$arr[2,12] = 0.7274271648586496
$arr[2,13] = 0.7274271648586496
$arr[3,0] = 20.002957
$arr[3,1] = 60.008871
$arr[3,2] = 0.7920860939997775
$arr[3,3] = 0.7920860939997775
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.06924866666666667
$arr[3,7] = 0.207746
$arr[3,8] = 0.0006844763088715736
$arr[3,9] = 0.0006844763088715734
$arr[3,10] = 1.385178101640667
$arr[3,11] = 12.466602914766
$arr[3,12] = 0.00054216416592947
$arr[3,13] = 0.0005421641659294699
$arr[4,0] = 1.047813333333333
$arr[4,1] = 3.14344
$arr[4,2] = 0.04149178396178559
$arr[4,3] = 0.04149178396178559
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 3.151158666666667
$arr[4,7] = 9.453476
$arr[4,8] = 0.03114707555614071
$arr[4,9] = 0.03114707555614071
$arr[4,10] = 3.301826066382222
$arr[4,11] = 29.71643459744
$arr[4,12] = 0.001292347730016803
$arr[4,13] = 0.001292347730016803
$arr[5,0] = 1.047813333333333
$arr[5,1] = 3.14344
$arr[5,2] = 0.04149178396178559
$arr[5,3] = 0.04149178396178559
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 5.038243666666667
$arr[5,7] = 15.114731
$arr[5,8] = 0.04979963650066307
$arr[5,9] = 0.04979963650066306
$arr[5,10] = 5.279138890515556
$arr[5,11] = 47.51225001464001
$arr[5,12] = 0.002066275759060964
$arr[5,13] = 0.002066275759060964
$arr[6,0] = 1.047813333333333
$arr[6,1] = 3.14344
$arr[6,2] = 0.04149178396178559
$arr[6,3] = 0.04149178396178559
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 92.91163899999999
$arr[6,7] = 278.734917
$arr[6,8] = 0.9183688116343246
$arr[6,9] = 0.9183688116343246
$arr[6,10] = 97.35405416605332
$arr[6,11] = 876.1864874944799
$arr[6,12] = 0.03810476032957316
$arr[6,13] = 0.03810476032957316
$arr[7,0] = 1.047813333333333
$arr[7,1] = 3.14344
$arr[7,2] = 0.04149178396178559
$arr[7,3] = 0.04149178396178559
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.06924866666666667
$arr[7,7] = 0.207746
$arr[7,8] = 0.0006844763088715736
$arr[7,9] = 0.0006844763088715734
$arr[7,10] = 0.07255967624888889
$arr[7,11] = 0.65303708624
$arr[7,12] = [double]"2.840014313465976E-05"
$arr[7,13] = [double]"2.840014313465975E-05"
$arr[8,0] = 2.298356333333333
$arr[8,1] = 6.895068999999999
$arr[8,2] = 0.09101134850660582
$arr[8,3] = 0.09101134850660582
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 3.151158666666667
$arr[8,7] = 9.453476
$arr[8,8] = 0.03114707555614071
$arr[8,9] = 0.03114707555614071
$arr[8,10] = 7.242485478871555
$arr[8,11] = 65.18236930984399
$arr[8,12] = 0.002834737348401505
$arr[8,13] = 0.002834737348401505
$arr[9,0] = 2.298356333333333
$arr[9,1] = 6.895068999999999
$arr[9,2] = 0.09101134850660582
$arr[9,3] = 0.09101134850660582
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 5.038243666666667
$arr[9,7] = 15.114731
$arr[9,8] = 0.04979963650066307
$arr[9,9] = 0.04979963650066306
$arr[9,10] = 11.57967924015989
$arr[9,11] = 104.217113161439
$arr[9,12] = 0.004532332073064135
$arr[9,13] = 0.004532332073064134
$arr[10,0] = 2.298356333333333
$arr[10,1] = 6.895068999999999
$arr[10,2] = 0.09101134850660582
$arr[10,3] = 0.09101134850660582
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 92.91163899999999
$arr[10,7] = 278.734917
$arr[10,8] = 0.9183688116343246
$arr[10,9] = 0.9183688116343246
$arr[10,10] = 213.5440539360303
$arr[10,11] = 1921.896485424273
$arr[10,12] = 0.08358198397324895
$arr[10,13] = 0.08358198397324895
$arr[11,0] = 2.298356333333333
$arr[11,1] = 6.895068999999999
$arr[11,2] = 0.09101134850660582
$arr[11,3] = 0.09101134850660582
$arr[11,4] = 2
$arr[11,5] = 0.6666666666666666
$arr[11,6] = 0.06924866666666667
$arr[11,7] = 0.207746
$arr[11,8] = 0.0006844763088715736
$arr[11,9] = 0.0006844763088715734
$arr[11,10] = 0.1591581116082222
$arr[11,11] = 1.432423004474
$arr[11,12] = [double]"6.229511189122595E-05"
$arr[11,13] = [double]"6.229511189122594E-05"
$arr[12,0] = 1.904387
$arr[12,1] = 5.713160999999999
$arr[12,2] = 0.07541077353183102
$arr[12,3] = 0.07541077353183102
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 3.151158666666667
$arr[12,7] = 9.453476
$arr[12,8] = 0.03114707555614071
$arr[12,9] = 0.03114707555614071
$arr[12,10] = 6.001025599737333
$arr[12,11] = 54.009230397636
$arr[12,12] = 0.002348825060942957
$arr[12,13] = 0.002348825060942957
$arr[13,0] = 1.904387
$arr[13,1] = 5.713160999999999
$arr[13,2] = 0.07541077353183102
$arr[13,3] = 0.07541077353183102
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 5.038243666666667
$arr[13,7] = 15.114731
$arr[13,8] = 0.04979963650066307
$arr[13,9] = 0.04979963650066306
$arr[13,10] = 9.594765741632333
$arr[13,11] = 86.35289167469101
$arr[13,12] = 0.003755429110119009
$arr[13,13] = 0.003755429110119008
$arr[14,0] = 1.904387
$arr[14,1] = 5.713160999999999
$arr[14,2] = 0.07541077353183102
$arr[14,3] = 0.07541077353183102
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 92.91163899999999
$arr[14,7] = 278.734917
$arr[14,8] = 0.9183688116343246
$arr[14,9] = 0.9183688116343246
$arr[14,10] = 176.939717460293
$arr[14,11] = 1592.457457142637
$arr[14,12] = 0.06925490247285283
$arr[14,13] = 0.06925490247285283
$arr[15,0] = 1.904387
$arr[15,1] = 5.713160999999999
$arr[15,2] = 0.07541077353183102
$arr[15,3] = 0.07541077353183102
$arr[15,4] = 2
$arr[15,5] = 0.6666666666666666
$arr[15,6] = 0.06924866666666667
$arr[15,7] = 0.207746
$arr[15,8] = 0.0006844763088715736
$arr[15,9] = 0.0006844763088715734
$arr[15,10] = 0.1318762605673333
$arr[15,11] = 1.186886345106
$arr[15,12] = [double]"5.161688791621785E-05"
$arr[15,13] = [double]"5.161688791621785E-05"
$ws.Range("G2:T17").Value = $arr
Write-Output "done"
